# Apply the "Updated cryptos list" data refresh to Sheet1.
# Cell values are set directly via Range.Value. Columns that contain
# numeric-looking text (Price) are forced to Text format first so Excel
# keeps the exact original formatting (e.g. "12.50", "2.90", "1.978.17")
# instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.905.91"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.311.00"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.59"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.75"
$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -1.17%  "

$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.39"
$ws.Range("E11").Value = "  +6.57%  "

$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.670.74"
$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.347.55"
$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.852.56"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.50"
$ws.Range("E19").Value = "  -1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  -1.33%  "

$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.69"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.65"
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("E24").Value = "  +3.73%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  -1.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.72"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  +15.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.15"
$ws.Range("E29").Value = "  -1.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.08"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.59"
$ws.Range("E31").Value = "  -1.44%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.98"
$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.73"
$ws.Range("E34").Value = "  +3.02%  "

$ws.Range("E35").Value = "  -6.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0699"
$ws.Range("E36").Value = "  +0.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -1.67%  "

$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("E40").Value = "  +0.75%  "

$ws.Range("E41").Value = "  -1.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.977.11"
$ws.Range("E42").Value = "  -1.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.58"
$ws.Range("E43").Value = "  +5.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.74"
$ws.Range("E44").Value = "  +4.31%  "

$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.77"
$ws.Range("E47").Value = "  -0.47%  "

# Row 48/49: RocketPoolETH and HuobiToken swap positions.
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.90"
$ws.Range("E48").Value = "  -0.13%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.536.77"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.39"
$ws.Range("E50").Value = "  -1.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.11"
